$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 4.75
$ws.Range("L2").Value = 5.5
$ws.Range("N2").Value = 7.5
$ws.Range("W2").Value = 5.5
$ws.Range("Z2").Value = 13
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 17
$ws.Range("AU2").Value = 9.5
$ws.Range("AW2").Value = 6.5
$ws.Range("AX2").Value = 29

# Row 3 updates
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 15
$ws.Range("Q3").Value = 1.73
$ws.Range("R3").Value = 2.08
